# Update the auto-updating "datetime" footer placeholder from 1/12/2018 to
# 1/15/2018 everywhere it appears: the slide master, every slide layout
# (custom layout), and the notes master.
#
# ppPlaceholderDate == 16; every "Date Placeholder" shape in this deck
# reports that PlaceholderFormat.Type, so use it instead of hard-coded
# shape indexes (they differ from layout to layout).

$newDate = "1/15/2018"

function Set-DatePlaceholderText($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$p = $ppt.ActivePresentation

# 1. Slide master.
$master = $p.SlideMaster
Set-DatePlaceholderText $master

# 2. Every slide layout off that master.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Set-DatePlaceholderText $layouts.Item($j)
}

# 3. Notes master.
Set-DatePlaceholderText $p.NotesMaster
